$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table column widths: 530 -> 529 dxa (col 1), 4710 -> 4711 dxa (col 2).
#    Word exposes column widths in points (1 pt = 20 dxa), so 530/20=26.5 ->
#    529/20=26.45, and 4710/20=235.5 -> 4711/20=235.55. Setting Columns(i).Width
#    updates every row uniformly (gridCol + every tcW in the column).
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 26.45
$t.Columns.Item(2).Width = 235.55

# ---------------------------------------------------------------------------
# 2) Merge the "4." + "1" runs (and the trailing ") " + "| ..." / ")" + "."
#    runs) into single runs, with no visible text change. Word's Find/Replace
#    leaves identical text alone, so we briefly swap in a placeholder string
#    first to force the engine to re-merge the run with unified formatting.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("предвзято (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $rng.End
$tmp = $d.Range($pos, $pos + 3)
$tmp.Text = "@@@@@"
$final = $d.Range($pos, $pos + 5)
$final.Text = "4.1"

$rng2 = $d.Content
$rng2.Find.Execute("4.3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos2 = $rng2.End
$tmp2 = $d.Range($pos2, $pos2 + 2)
$tmp2.Text = "@@"
$final2 = $d.Range($pos2, $pos2 + 2)
$final2.Text = ")."

# ---------------------------------------------------------------------------
# 3) Shorten "Отношение к проблемной ситуации: негативно-нейтральное ..." --
#    drop everything after "негативно-нейтральное" up to (but excluding) the
#    trailing period, which stays as its own (non-bold) run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " со стороны именно чистого научного биоинформатического сообщества России; позитивное со стороны организаций, являющихся потенциальными «конкурентами» Института биоинформатики",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Shorten "... особое доверие (7.1) | улучшение дел ... крах (7.2)." down
#    to "... особое доверие (7)." -- keep the bold/italic/red "7" marker run
#    and the closing ")." run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "7.1) | улучшение дел у потенциальных «конкурентов», т.к. такой сильный игрок на научном рынке, как Институт биоинформатики, постепенно терпит свой крах (7.2",
    $true, $false, $false, $false, $false, $true, 1, $false, "7", 2) | Out-Null
